# Natmi following Dr Hou advice
# Rebuild the LR-pair table with a third "ECs" sending/target cluster category.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$clusters = @("ECs", "FAPs", "sCs")

$rowData = @{
    "ECs|ECs"   = @(2, 0.6666666666666666, 3.675031333333333, 11.025094, 0.2032371147293133, 0.2032371147293133, 3, 1, 114.155417, 342.466251, 0.6835107367845005, 0.6835107367845005, 419.5247343447326, 3775.722609102594, 0.138914750030589, 0.138914750030589)
    "ECs|FAPs"  = @(2, 0.6666666666666666, 3.675031333333333, 11.025094, 0.2032371147293133, 0.2032371147293133, 3, 1, 35.924535, 107.773605, 0.2150997826628812, 0.2150997826628812, 132.02379176043, 1188.21412584387, 0.04371625920730634, 0.04371625920730635)
    "ECs|sCs"   = @(2, 0.6666666666666666, 3.675031333333333, 11.025094, 0.2032371147293133, 0.2032371147293133, 3, 1, 16.93339666666667, 50.80019, 0.1013894805526183, 0.1013894805526183, 62.23076332976222, 560.0768699678599, 0.02060610549141798, 0.02060610549141798)
    "FAPs|ECs"  = @(3, 1, 10.108494, 30.325482, 0.5590213983169419, 0.5590213983169419, 3, 1, 114.155417, 342.466251, 0.6835107367845005, 0.6835107367845005, 1153.939347811998, 10385.45413030798, 0.3820971278419146, 0.3820971278419146)
    "FAPs|FAPs" = @(3, 1, 10.108494, 30.325482, 0.5590213983169419, 0.5590213983169419, 3, 1, 35.924535, 107.773605, 0.2150997826628812, 0.2150997826628812, 363.14294650029, 3268.28651850261, 0.1202453812818741, 0.1202453812818741)
    "FAPs|sCs"  = @(3, 1, 10.108494, 30.325482, 0.5590213983169419, 0.5590213983169419, 3, 1, 16.93339666666667, 50.80019, 0.1013894805526183, 0.1013894805526183, 171.17113860462, 1540.54024744158, 0.05667888919315309, 0.05667888919315309)
    "sCs|ECs"   = @(3, 1, 4.298956, 12.896868, 0.2377414869537448, 0.2377414869537448, 3, 1, 114.155417, 342.466251, 0.6835107367845005, 0.6835107367845005, 490.7491148446521, 4416.742033601869, 0.1624988589119968, 0.1624988589119968)
    "sCs|FAPs"  = @(3, 1, 4.298956, 12.896868, 0.2377414869537448, 0.2377414869537448, 3, 1, 35.924535, 107.773605, 0.2150997826628812, 0.2150997826628812, 154.43799528546, 1389.94195756914, 0.0511381421737007, 0.05113814217370071)
    "sCs|sCs"   = @(3, 1, 4.298956, 12.896868, 0.2377414869537448, 0.2377414869537448, 3, 1, 16.93339666666667, 50.80019, 0.1013894805526183, 0.1013894805526183, 72.79592720054667, 655.16334480492, 0.02410448586804727, 0.02410448586804727)
}

$r = 2
foreach ($sender in $clusters) {
    foreach ($target in $clusters) {
        $ws.Cells.Item($r, 1).Value = $sender
        $ws.Cells.Item($r, 2).Value = "Inhba"
        $ws.Cells.Item($r, 3).Value = "Eng"
        $ws.Cells.Item($r, 4).Value = $target

        $vals = $rowData["$sender|$target"]
        for ($i = 0; $i -lt $vals.Length; $i++) {
            $ws.Cells.Item($r, 5 + $i).Value = $vals[$i]
        }
        $r++
    }
}
